# Update dashboards - 2026-03-01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22 (PPI-FD / PPIFIS, M/M % Delta) ---
$ws.Range("N22").Value = 46023
$ws.Range("N22").NumberFormat = "yyyy-mm-dd"
$ws.Range("N22").Interior.Color = 65535
$ws.Range("Q22").Value = 0.00480699391202144
$ws.Range("R22").Value = 0.004303799147197918
$ws.Range("S22").Value = 0.00248632190555953
$ws.Range("T22").Value = 0.001184731605045064
$ws.Range("U22").Value = 0.006147582151921682

# --- Row 23 (PPIFIS, Y/Y % Delta) ---
$ws.Range("N23").Value = 46023
$ws.Range("N23").NumberFormat = "yyyy-mm-dd"
$ws.Range("N23").Interior.Color = 65535
$ws.Range("Q23").Value = 0.02837641493495531
$ws.Range("R23").Value = 0.02988058645921175
$ws.Range("S23").Value = 0.02970357876911097
$ws.Range("T23").Value = 0.02822398731321852
$ws.Range("U23").Value = 0.03046576545064236

# --- Row 29 (T5YIFR) ---
$ws.Range("N29").Value = 46080
$ws.Range("Q29").Value = 2.1
$ws.Range("R29").Value = 2.13
$ws.Range("S29").Value = 2.14
$ws.Range("T29").Value = 2.12
$ws.Range("U29").Value = 2.12

# --- Row 30 (T10YIE) ---
$ws.Range("N30").Value = 46080
$ws.Range("Q30").Value = 2.25
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.28
$ws.Range("T30").Value = 2.26
$ws.Range("U30").Value = 2.26

# --- Row 47 (DFF) ---
$ws.Range("N47").Value = 46079

# --- Row 48 (DGS2) ---
$ws.Range("N48").Value = 46079
$ws.Range("Q48").Value = 3.42
$ws.Range("R48").Value = 3.45
$ws.Range("S48").Value = 3.43
$ws.Range("T48").Value = 3.43
$ws.Range("U48").Value = 3.48

# --- Row 49 (DGS5) ---
$ws.Range("N49").Value = 46079
$ws.Range("Q49").Value = 3.57
$ws.Range("R49").Value = 3.61
$ws.Range("S49").Value = 3.61
$ws.Range("T49").Value = 3.59
$ws.Range("U49").Value = 3.65

# --- Row 50 (DGS10) ---
$ws.Range("N50").Value = 46079
$ws.Range("Q50").Value = 4.02
$ws.Range("R50").Value = 4.05
$ws.Range("S50").Value = 4.04
$ws.Range("T50").Value = 4.03
$ws.Range("U50").Value = 4.08
